$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (prevents Excel from
# auto-coercing numeric-looking strings like "294.97" into floating point
# numbers, and avoids leaving any lingering NumberFormat/style on the cell).
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "40.144.46"
Set-TextValue "E2" "  +0.49%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.223.09"
Set-TextValue "E3" "  +0.52%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "294.97"
Set-TextValue "E5" "  +1.77%  "

# Row 6 - Solana
Set-TextValue "D6" "87.82"
Set-TextValue "E6" "  +0.74%  "

# Row 7 - XRP
Set-TextValue "E7" "  +0.76%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.473"
Set-TextValue "E9" "  +0.74%  "

# Row 10 - was OKB, now Avalanche
Set-TextValue "B10" "Avalanche"
Set-TextValue "C10" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D10" "30.87"
Set-TextValue "E10" "  +1.87%  "

# Row 11 - was Avalanche, now OKB
Set-TextValue "B11" "OKB"
Set-TextValue "C11" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D11" "51.57"
Set-TextValue "E11" "  +8.13%  "

# Row 12 - Dogecoin
Set-TextValue "E12" "  +1.27%  "

# Row 13 - TRON
Set-TextValue "E13" "  +3.95%  "

# Row 14 - Polkadot
Set-TextValue "E14" "  -0.35%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.566.85"
Set-TextValue "E15" "  +0.44%  "

# Row 16 - Chainlink
Set-TextValue "D16" "13.91"
Set-TextValue "E16" "  +0.04%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.237.20"
Set-TextValue "E17" "  +1.92%  "

# Row 18 - Polygon
Set-TextValue "D18" "0.737"
Set-TextValue "E18" "  +1.58%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "40.073.91"
Set-TextValue "E19" "  +0.45%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.0₃0889"
Set-TextValue "E20" "  +1.02%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "E21" "  -1.66%  "

# Row 22 - Uniswap
Set-TextValue "E22" "  +0.26%  "

# Row 23 - Litecoin
Set-TextValue "D23" "65.71"
Set-TextValue "E23" "  +0.37%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "235.74"
Set-TextValue "E24" "  +0.17%  "

# Row 25 - Dai
Set-TextValue "E25" "  +0.03%  "

# Row 26 - PancakeSwap
Set-TextValue "E26" "  +1.65%  "

# Row 27 - ImmutableX
Set-TextValue "E27" "  +0.51%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "23.26"
Set-TextValue "E28" "  +3.69%  "

# Row 29 - Cosmos
Set-TextValue "E29" "  +1.89%  "

# Row 30 - Toncoin
Set-TextValue "E30" "  -4.88%  "

# Row 31 - Monero
Set-TextValue "D31" "161.75"
Set-TextValue "E31" "  +4.00%  "

# Row 32 - InjectiveProtocol
Set-TextValue "D32" "31.81"
Set-TextValue "E32" "  +0.63%  "

# Row 33 - FirstDigitalUSD
Set-TextValue "E33" "  -0.02%  "

# Row 34 - LidoDAOToken
Set-TextValue "E34" "  +8.24%  "

# Row 35 - Filecoin
Set-TextValue "E35" "  +1.39%  "

# Row 36 - Hedera
Set-TextValue "E36" "  +0.30%  "

# Row 37 - WEMIXToken
Set-TextValue "E37" "  -1.31%  "

# Row 38 - Stellar
Set-TextValue "E38" "  +1.90%  "

# Row 39 - ARBITRUM
Set-TextValue "E39" "  +4.84%  "

# Row 40 - Kaspa
Set-TextValue "E40" "  +1.63%  "

# Row 41 - Celestia
Set-TextValue "D41" "15.68"
Set-TextValue "E41" "  -0.45%  "

# Row 42 - Maker
Set-TextValue "D42" "2.081.19"
Set-TextValue "E42" "  -0.93%  "

# Row 43 - RenderToken
Set-TextValue "D43" "3.76"
Set-TextValue "E43" "  -1.36%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "19.71"
Set-TextValue "E44" "  +12.91%  "

# Row 45 - VeChain
Set-TextValue "E45" "  +1.72%  "

# Row 46 - FraxShare
Set-TextValue "D46" "9.94"
Set-TextValue "E46" "  +0.08%  "

# Row 47 - NEARProtocol
Set-TextValue "E47" "  +5.26%  "

# Row 48 - ApeXProtocol
Set-TextValue "E48" "  -10.79%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "2.439.07"
Set-TextValue "E49" "  +0.45%  "

# Row 50 - TrustWalletToken
Set-TextValue "E50" "  +4.18%  "

# Row 51 - Stacks
Set-TextValue "E51" "  +2.25%  "
